# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.760.84"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "2.282.48"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "124.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.641"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.40%  "

$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.37%  "

$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.906"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.84%  "

$ws.Range("D16").Value = "2.625.88"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "2.276.10"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "43.719.98"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.55%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.23%  "

$ws.Range("E27").Value = "  +1.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.69%  "

$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0926"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.93%  "

$ws.Range("E36").Value = "  +2.10%  "

$ws.Range("E37").Value = "  +5.09%  "

$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.69%  "

$ws.Range("E43").Value = "  -0.72%  "

$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +37.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.72%  "

$ws.Range("E49").Value = "  +0.82%  "

$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
